# Updates the hme_transformed report: the report date moves forward one day
# (2026-02-04 -> 2026-02-05) and the per-store/day-part metrics in columns
# D:O are refreshed with the newly downloaded Gmail data. Columns B (store)
# and C (time_measure) are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 30

# --- Column A: new report date for every data row ---------------------
$newDate = Get-Date -Year 2026 -Month 2 -Day 5 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
$dateRange = $ws.Range("A$firstRow`:A$lastRow")
$numDateRows = $lastRow - $firstRow + 1
$dateValues = New-Object 'object[,]' $numDateRows, 1
for ($i = 0; $i -lt $numDateRows; $i++) {
    $dateValues[$i, 0] = $newDate
}
$dateRange.Value = $dateValues

# --- Columns D:O: refreshed metrics ------------------------------------
$data = @(
    @(208,28,2,28,2,0,0,0,0,44,47,119),
    @(103,22,2,22,2,0,0,0,0,51,29,103),
    @(56,22,2,22,2,0,0,0,0,71,14,108),
    @(12,28,3,28,3,0,0,0,0,76,32,137),
    @(5,17,3,17,3,0,0,0,0,53,8,80),
    @(267,28,2,28,2,0,0,0,0,43,55,127),
    @(115,27,2,27,2,0,0,0,0,45,35,108),
    @(93,28,3,28,3,0,0,0,0,66,40,135),
    @(17,31,2,31,2,0,0,0,0,75,16,120),
    @(248,26,3,26,3,0,0,0,0,48,70,143),
    @(99,37,9,37,9,0,0,0,0,74,45,156),
    @(75,59,12,59,12,0,0,0,0,88,34,181),
    @(19,51,2,51,2,0,0,0,0,115,26,193),
    @(210,30,6,30,6,0,0,0,0,52,91,172),
    @(78,42,7,42,7,0,0,0,0,80,40,161),
    @(45,58,7,58,7,0,0,0,0,92,32,181),
    @(9,31,3,31,3,0,0,0,0,90,18,140),
    @(167,33,3,33,3,0,0,0,0,51,40,95),
    @(50,37,3,37,3,0,0,0,0,88,28,129),
    @(39,48,5,48,5,0,0,0,0,126,36,183),
    @(8,61,6,61,6,0,0,0,0,210,11,273),
    @(225,26,2,0,0,26,2,0,0,48,50,123),
    @(71,34,3,0,0,34,3,0,0,97,98,230),
    @(55,36,3,0,0,36,3,0,0,75,28,141),
    @(12,27,2,0,0,27,2,0,0,68,9,105),
    @(223,27,2,27,2,0,0,0,0,46,52,138),
    @(98,28,2,28,2,0,0,0,0,58,34,128),
    @(75,43,3,43,3,0,0,0,0,97,56,205),
    @(14,23,3,23,3,0,0,0,0,107,20,151)
)

$numRows = $data.Count
$numCols = $data[0].Count
$values = New-Object 'object[,]' $numRows, $numCols
for ($r = 0; $r -lt $numRows; $r++) {
    for ($c = 0; $c -lt $numCols; $c++) {
        $values[$r, $c] = $data[$r][$c]
    }
}

$lastDataRow = $firstRow + $numRows - 1
$metricsRange = $ws.Range("D$firstRow`:O$lastDataRow")
$metricsRange.Value = $values

$wb.Save()
